$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6 (shifts Doveton..Southbank down by one row)
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the Camberwell entry
$ws.Range("A6").Value = "Camberwell"
$ws.Range("B6").Value = "Tao Dumplings  1 Evans Place, Camberwell VIC 3124"
$ws.Range("C6").Value = "29/12/20 12:30pm-1:30pm"
$ws.Range("D6").Value = "Case ate at restaurant"

# Fix the exposure period time for the Melbourne Central Lion Hotel row,
# which is now row 11 after the insertion above (10:30pm -> 10:00pm)
$ws.Range("C11").Value = "28/12/2020 10:00pm-12.00am"
